$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Chapter 9: 2 Knights Def (II) - mark as Done* with note
$ws.Range("E5").Value = "Done*"
$ws.Range("E5").Font.Bold = $true
$ws.Range("E5").Font.Italic = $false
$ws.Range("F5").Value = "All lines followed at least to move 8"

# Row 6: Chapter 10: 2 Knights Def (III) - add new training file, mark Done*
$ws.Range("D6").Value = "E:\Chess\Database\Openings\Open-Games-e4-e5\Two-Knights-Main-Line-7...Bc5.pgn"
$ws.Range("E6").Value = "Done*"
$ws.Range("E6").Font.Bold = $true
$ws.Range("F6").Value = "All lines followed at least to move 8"

# Move selection to D7 to reflect final cursor position
$ws.Range("D7").Select()
